# Add data for 2022-11-29: the "through" date in the report moves from
# 11-20 to 11-21, which bumps the November figures (and therefore the
# Total row) for every year column except 2015 and 2018.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Workbook/sheet title: "Through 2022-11-20" -> "Through 2022-11-21"
$ws.Name = "Through 2022-11-21"

# Row label in column A: "November (through 11-20)" -> "November (through 11-21)"
$ws.Range("A12").Value = "November (through 11-21)"

# Updated November row (row 12) and Total row (row 13) values, per year
# column. Columns B (2015) and E (2018) are unchanged.
$updates = @{
    "C" = @{ Nov = 50;  Total = 536 }   # 2016
    "D" = @{ Nov = 85;  Total = 795 }   # 2017
    "F" = @{ Nov = 32;  Total = 514 }   # 2019
    "G" = @{ Nov = 143; Total = 1200 }  # 2020
    "H" = @{ Nov = 149; Total = 1590 }  # 2021
    "I" = @{ Nov = 81;  Total = 1478 }  # 2022
}

foreach ($col in $updates.Keys) {
    $ws.Range("$col" + "12").Value = $updates[$col].Nov
    $ws.Range("$col" + "13").Value = $updates[$col].Total
}
